$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rankings: append new row 19 (2020, 38)
# ------------------------------------------------------------------
$wsRankings = $wb.Worksheets.Item("Rankings")
$wsRankings.Cells.Item(19,1).Value = 2020
$wsRankings.Cells.Item(19,2).Value = 38

# ------------------------------------------------------------------
# 2. EngageLearn: insert new row 2 (2021, 1229, 11722, 21946)
# ------------------------------------------------------------------
$wsEngageLearn = $wb.Worksheets.Item("EngageLearn")
$wsEngageLearn.Rows.Item(2).Insert()
$wsEngageLearn.Cells.Item(2,1).Value = 2021
$wsEngageLearn.Cells.Item(2,2).Value = 1229
$wsEngageLearn.Cells.Item(2,3).Value = 11722
$wsEngageLearn.Cells.Item(2,4).Value = 21946

# ------------------------------------------------------------------
# 3. Collections: insert new row 2 (2021, 2857879, 3114927, 39627)
# ------------------------------------------------------------------
$wsCollections = $wb.Worksheets.Item("Collections")
$wsCollections.Rows.Item(2).Insert()
$wsCollections.Cells.Item(2,1).Value = 2021
$wsCollections.Cells.Item(2,2).Value = 2857879
$wsCollections.Cells.Item(2,3).Value = 3114927
$wsCollections.Cells.Item(2,4).Value = 39627

# ------------------------------------------------------------------
# 4. CollectionUseDelivery: insert new row 2
#    (2021, 17261, 3071624, 690754, 7845811, 1218407, 11962, 13457)
# ------------------------------------------------------------------
$wsCUD = $wb.Worksheets.Item("CollectionUseDelivery")
$wsCUD.Rows.Item(2).Insert()
$wsCUD.Cells.Item(2,1).Value = 2021
$wsCUD.Cells.Item(2,2).Value = 17261
$wsCUD.Cells.Item(2,3).Value = 3071624
$wsCUD.Cells.Item(2,4).Value = 690754
$wsCUD.Cells.Item(2,5).Value = 7845811
$wsCUD.Cells.Item(2,6).Value = 1218407
$wsCUD.Cells.Item(2,7).Value = 11962
$wsCUD.Cells.Item(2,8).Value = 13457

# ------------------------------------------------------------------
# 5. Expenditures: insert new row 2 (2021, 19.67, 76.51, 3.8)
# ------------------------------------------------------------------
$wsExpenditures = $wb.Worksheets.Item("Expenditures")
$wsExpenditures.Rows.Item(2).Insert()
$wsExpenditures.Cells.Item(2,1).Value = 2021
$wsExpenditures.Cells.Item(2,2).Value = 19.67
$wsExpenditures.Cells.Item(2,3).Value = 76.51
$wsExpenditures.Cells.Item(2,4).Value = 3.8

# ------------------------------------------------------------------
# 6. SpacesStaff: insert new row 2 (2021, 10, 311555, 157)
# ------------------------------------------------------------------
$wsSpacesStaff = $wb.Worksheets.Item("SpacesStaff")
$wsSpacesStaff.Rows.Item(2).Insert()
$wsSpacesStaff.Cells.Item(2,1).Value = 2021
$wsSpacesStaff.Cells.Item(2,2).Value = 10
$wsSpacesStaff.Cells.Item(2,3).Value = 311555
$wsSpacesStaff.Cells.Item(2,4).Value = 157

# ------------------------------------------------------------------
# 7. Visitors: insert new row 2 (2021, 251916, 10894129)
#    old row2 (now row3) gets a custom row height of 14.25
# ------------------------------------------------------------------
$wsVisitors = $wb.Worksheets.Item("Visitors")
$wsVisitors.Rows.Item(2).Insert()
$wsVisitors.Cells.Item(2,1).Value = 2021
$wsVisitors.Cells.Item(2,2).Value = 251916
$wsVisitors.Cells.Item(2,3).Value = 10894129
$wsVisitors.Rows.Item(3).RowHeight = 14.25

# ------------------------------------------------------------------
# 8. SocialMedia: insert new row 2 (2021, 2573, 1082, 905)
# ------------------------------------------------------------------
$wsSocialMedia = $wb.Worksheets.Item("SocialMedia")
$wsSocialMedia.Rows.Item(2).Insert()
$wsSocialMedia.Cells.Item(2,1).Value = 2021
$wsSocialMedia.Cells.Item(2,2).Value = 2573
$wsSocialMedia.Cells.Item(2,3).Value = 1082
$wsSocialMedia.Cells.Item(2,4).Value = 905

# ------------------------------------------------------------------
# 9. New sheet "DigitalCollection" added after SocialMedia
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDigital = $wb.Worksheets.Add($null, $lastSheet)
$wsDigital.Name = "DigitalCollection"

$wsDigital.Cells.Item(1,1).Value = "Notable Digital Collections"
$wsDigital.Cells.Item(1,2).Value = "Number of Digital Items"

$wsDigital.Cells.Item(2,1).Value = "Television News Archive"
$wsDigital.Cells.Item(2,2).Value = 1293716

$wsDigital.Cells.Item(3,1).Value = "History of Art Image Repository"
$wsDigital.Cells.Item(3,2).Value = 106259

$wsDigital.Cells.Item(4,1).Value = "Institutional Respository"
$wsDigital.Cells.Item(4,2).Value = 9687

$wsDigital.Cells.Item(5,1).Value = "Fine Arts Gallery"
$wsDigital.Cells.Item(5,2).Value = 7113

$wsDigital.Cells.Item(6,1).Value = "Art in the Christian Tradition"
$wsDigital.Cells.Item(6,2).Value = 6174

$wsDigital.Cells.Item(7,1).Value = "The N$([char]0x00FC)rnberg Krupp Trial Papers of Judge Hu C. Anderson"
$wsDigital.Cells.Item(7,2).Value = 6052

$wsDigital.Cells.Item(8,1).Value = "ETD, Electronic Theses & Dissertations"
$wsDigital.Cells.Item(8,2).Value = 5535

$wsDigital.Cells.Item(9,1).Value = "History of Medicine Collection"
$wsDigital.Cells.Item(9,2).Value = 5510

$wsDigital.Cells.Item(10,1).Value = "Blair Performance Archive"
$wsDigital.Cells.Item(10,2).Value = 4955

$wsDigital.Cells.Item(11,1).Value = "Helguera Collection of Colombiana"
$wsDigital.Cells.Item(11,2).Value = 3093

$wsDigital.Cells.Item(12,1).Value = "Manual Zapata Olivella"
$wsDigital.Cells.Item(12,2).Value = 2573

$wsDigital.Cells.Item(13,1).Value = "Global Music Archive"
$wsDigital.Cells.Item(13,2).Value = 1805

$wsDigital.Cells.Item(14,1).Value = "Revised Common Lectionary"
$wsDigital.Cells.Item(14,2).Value = 249

$wsDigital.Columns.Item(1).ColumnWidth = 52.86
$wsDigital.Columns.Item(2).ColumnWidth = 22.57

$wsDigital.Range("B2").Select()

# ------------------------------------------------------------------
# 10. Fix up selections on each sheet to match the final saved state
# ------------------------------------------------------------------
$wsEngageLearn.Range("E2").Select()
$wsCollections.Range("E2").Select()
$wsCUD.Range("I2").Select()
$wsExpenditures.Range("E2").Select()
$wsSpacesStaff.Range("E2").Select()
$wsVisitors.Range("D11").Select()
$wsSocialMedia.Range("E2").Select()

# ------------------------------------------------------------------
# 11. Final active sheet/selection: Rankings!C19
# ------------------------------------------------------------------
$wsRankings.Activate()
$wsRankings.Range("C19").Select()
